# Converts the old-style heading/byline paragraphs into the pandoc-flavoured
# "title block" layout: a Title-styled heading and an Authors-styled byline,
# each split word-by-word into its own run (the way pandoc's docx writer
# renders its title-block Inlines).
#
# w:document.xml (before):
#   <w:bookmarkStart .../>
#   <w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr>
#     <w:r><w:t>Chicago Readers</w:t></w:r></w:p>
#   <w:bookmarkEnd .../>
#   <w:p><w:r><w:rPr><w:b/></w:rPr><w:t>By Dorothy Day</w:t></w:r></w:p>
#
# w:document.xml (after):
#   <w:p><w:pPr><w:pStyle w:val="Title"/></w:pPr>
#     <w:r><w:t>Chicago</w:t></w:r><w:r><w:t> </w:t></w:r><w:r><w:t>Readers</w:t></w:r></w:p>
#   <w:p><w:pPr><w:pStyle w:val="Authors"/></w:pPr>
#     <w:r><w:t>Dorothy</w:t></w:r><w:r><w:t> </w:t></w:r><w:r><w:t>Day</w:t></w:r></w:p>

$d = $word.ActiveDocument

function Set-ParagraphAsWordRuns($paraIndex, $styleName, $words) {
    # Replace the whole paragraph (minus its trailing pilcrim) with fresh
    # runs - one per element of $words - carrying no direct/inherited
    # character formatting and using paragraph style $styleName.
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range

    $innerXml = ""
    foreach ($w in $words) {
        $escaped = $w.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
        $innerXml += '<w:r><w:t xml:space="preserve">' + $escaped + '</w:t></w:r>'
    }

    $packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p><w:pPr><w:pStyle w:val="' + $styleName + '"/></w:pPr>' + $innerXml + '</w:p></w:body>' +
        '</w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($packageXml)
}

Set-ParagraphAsWordRuns 1 "Title" @("Chicago", " ", "Readers")
Set-ParagraphAsWordRuns 2 "Authors" @("Dorothy", " ", "Day")
